$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo "2 Hour" -> "2 Hours" on existing rows 3 and 4 ---
$ws.Range("A3").Value = "2 Hours"
$ws.Range("A4").Value = "2 Hours"

# --- Row 5: turn the lone date/text row into a full 4-column entry ---
# The existing row 5 carries a leftover "ht=30" row height from when it
# only held a date cell. Stash a copy of its (only) style-2 formatting in
# a scratch cell, delete the row so the rebuilt one picks up a plain,
# non-custom height, then recreate the cells from scratch.
$ws.Range("B5").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(5).Delete() | Out-Null

# A5 needs the "wrap/center" style (same as A3/A4 etc.)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null

# B5 needs the "wrap/center + centered date" style, restored from scratch cell
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# C5 needs the "wrap/center" style
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null

# D5 needs the "wrap/center" style too
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null

# Now set the row 5/6 values in the exact order that introduces new shared
# strings in the same order as the target workbook.
$ws.Range("C5").Value = "Draw Snow"
$ws.Range("D5").Value = "Created particle system for display the snow"
$ws.Range("B5").Value = "19/08/2024 1pm - 4pm"

# --- Row 6 (brand new row) ---
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null

$ws.Range("C6").Value = "Add Features"
$ws.Range("B6").Value = "20/08/2024 3pm - 5pm"

$ws.Range("A5").Value = "3 Hours"

$ws.Range("A6").Value = "2 Hours"
$ws.Range("D6").Value = " Fixed scaling system and the snow disappearing.  Display diagnostic data to the screen"

# --- Row 7 (brand new row, only B/C/D columns used) ---
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null

$ws.Range("B7").Value = "21/08/2024 9am - "
$ws.Range("D7").Value = "Added the ability for the snowman to jumping following a parabolic curve based on the time."
$ws.Range("C7").Value = "Add Features"

# --- Row heights to match the filled-in content ---
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 30

# --- Selection / active cell ---
$ws.Range("H6").Select() | Out-Null
